$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Recorded By" is column G (7). Walk every data row and, where the
# cell holds a comma-separated list of recorders, rotate the list left
# by one position (first entry moves to the end) - matching the
# reordering seen in the source data - except for the single list
# "dnasr281@gmail.com, System" which is left as-is.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq '') {
        continue
    }

    $parts = $val -split ', '
    if ($parts.Count -le 1) {
        continue
    }

    if ($parts.Count -eq 2 -and $parts[0] -eq 'dnasr281@gmail.com' -and $parts[1] -eq 'System') {
        continue
    }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ', '
    $cell.Value2 = $rotated
}
